# Update column 4 ("Dysmenorrhea.plus.Bladder.Pain") values in Table 1
# of the document per the target diff (table 4 + minor edits).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$editCount = 0

$cell = $t.Cell(2, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '76') {
    $rng.Text = '75'
    $editCount = $editCount + 1
} else {
    throw "Row 2 col 4: expected '76' but found '$($rng.Text)'"
}

$cell = $t.Cell(4, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '76 (100.0)') {
    $rng.Text = '75 (100.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 4 col 4: expected '76 (100.0)' but found '$($rng.Text)'"
}

$cell = $t.Cell(6, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '2 (2.6)') {
    $rng.Text = '2 (2.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 6 col 4: expected '2 (2.6)' but found '$($rng.Text)'"
}

$cell = $t.Cell(8, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '73 (96.1)') {
    $rng.Text = '72 (96.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 8 col 4: expected '73 (96.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '14 (18.4)') {
    $rng.Text = '14 (18.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 9 col 4: expected '14 (18.4)' but found '$($rng.Text)'"
}

$cell = $t.Cell(10, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '17 (22.4)') {
    $rng.Text = '16 (21.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 10 col 4: expected '17 (22.4)' but found '$($rng.Text)'"
}

$cell = $t.Cell(11, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '2 (2.6)') {
    $rng.Text = '2 (2.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 11 col 4: expected '2 (2.6)' but found '$($rng.Text)'"
}

$cell = $t.Cell(14, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '38 (50.0)') {
    $rng.Text = '38 (50.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 14 col 4: expected '38 (50.0)' but found '$($rng.Text)'"
}

$cell = $t.Cell(15, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '16 (21.1)') {
    $rng.Text = '16 (21.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 15 col 4: expected '16 (21.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(16, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '59 (77.6)') {
    $rng.Text = '58 (77.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 16 col 4: expected '59 (77.6)' but found '$($rng.Text)'"
}

$cell = $t.Cell(19, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '45 (59.2)') {
    $rng.Text = '45 (60.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 19 col 4: expected '45 (59.2)' but found '$($rng.Text)'"
}

$cell = $t.Cell(21, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '7 (9.2)') {
    $rng.Text = '6 (8.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 21 col 4: expected '7 (9.2)' but found '$($rng.Text)'"
}

$cell = $t.Cell(22, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '16 (21.1)') {
    $rng.Text = '16 (21.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 22 col 4: expected '16 (21.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(23, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '7 (9.2)') {
    $rng.Text = '7 (9.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 23 col 4: expected '7 (9.2)' but found '$($rng.Text)'"
}

$cell = $t.Cell(25, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '13 (17.1)') {
    $rng.Text = '12 (16.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 25 col 4: expected '13 (17.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(26, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '11 (14.5)') {
    $rng.Text = '11 (14.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 26 col 4: expected '11 (14.5)' but found '$($rng.Text)'"
}

$cell = $t.Cell(27, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '6 (7.9)') {
    $rng.Text = '6 (8.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 27 col 4: expected '6 (7.9)' but found '$($rng.Text)'"
}

$cell = $t.Cell(28, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '34 (44.7)') {
    $rng.Text = '34 (45.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 28 col 4: expected '34 (44.7)' but found '$($rng.Text)'"
}

$cell = $t.Cell(29, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '12 (15.8)') {
    $rng.Text = '12 (16.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 29 col 4: expected '12 (15.8)' but found '$($rng.Text)'"
}

$cell = $t.Cell(30, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '10 (13.2)') {
    $rng.Text = '10 (13.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 30 col 4: expected '10 (13.2)' but found '$($rng.Text)'"
}

$cell = $t.Cell(31, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '5 (6.6)') {
    $rng.Text = '5 (6.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 31 col 4: expected '5 (6.6)' but found '$($rng.Text)'"
}

$cell = $t.Cell(32, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '16 (21.1)') {
    $rng.Text = '16 (21.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 32 col 4: expected '16 (21.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(33, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '13 (17.1)') {
    $rng.Text = '13 (17.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 33 col 4: expected '13 (17.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(34, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '16 (21.1)') {
    $rng.Text = '16 (21.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 34 col 4: expected '16 (21.1)' but found '$($rng.Text)'"
}

$cell = $t.Cell(35, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '11 (14.5)') {
    $rng.Text = '11 (14.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 35 col 4: expected '11 (14.5)' but found '$($rng.Text)'"
}

$cell = $t.Cell(36, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '5 (6.6)') {
    $rng.Text = '4 (5.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 36 col 4: expected '5 (6.6)' but found '$($rng.Text)'"
}

$cell = $t.Cell(37, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '0.13 (0.50)') {
    $rng.Text = '0.11 (0.45)'
    $editCount = $editCount + 1
} else {
    throw "Row 37 col 4: expected '0.13 (0.50)' but found '$($rng.Text)'"
}

$cell = $t.Cell(38, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '0.05 (0.32)') {
    $rng.Text = '0.03 (0.23)'
    $editCount = $editCount + 1
} else {
    throw "Row 38 col 4: expected '0.05 (0.32)' but found '$($rng.Text)'"
}

$cell = $t.Cell(39, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '69 (90.8)') {
    $rng.Text = '69 (92.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 39 col 4: expected '69 (90.8)' but found '$($rng.Text)'"
}

$cell = $t.Cell(40, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '7 (9.2)') {
    $rng.Text = '6 (8.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 40 col 4: expected '7 (9.2)' but found '$($rng.Text)'"
}

$cell = $t.Cell(41, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '39 (51.3)') {
    $rng.Text = '39 (52.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 41 col 4: expected '39 (51.3)' but found '$($rng.Text)'"
}

$cell = $t.Cell(42, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '37 (48.7)') {
    $rng.Text = '36 (48.0)'
    $editCount = $editCount + 1
} else {
    throw "Row 42 col 4: expected '37 (48.7)' but found '$($rng.Text)'"
}

$cell = $t.Cell(43, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '74 (97.4)') {
    $rng.Text = '73 (97.3)'
    $editCount = $editCount + 1
} else {
    throw "Row 43 col 4: expected '74 (97.4)' but found '$($rng.Text)'"
}

$cell = $t.Cell(44, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
if ($rng.Text -eq '2 (2.6)') {
    $rng.Text = '2 (2.7)'
    $editCount = $editCount + 1
} else {
    throw "Row 44 col 4: expected '2 (2.6)' but found '$($rng.Text)'"
}

Write-Host "Applied $editCount cell edits."